# Scheduled market-data refresh: updates currentAveragePrice* / LeveProfit*
# columns (H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets with the
# latest pulled values.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 654.38464
$ws.Range("J17").Value = 678.04
$ws.Range("L17").Value = 2034.12
$ws.Range("N17").Value = -2370.12
$ws.Range("H18").Value = 4359.0835
$ws.Range("I18").Value = 5933.7144
$ws.Range("J18").Value = 2154.6
$ws.Range("K18").Value = 5933.7144
$ws.Range("L18").Value = 2154.6
$ws.Range("M18").Value = -5649.7144
$ws.Range("N18").Value = -2722.6
$ws.Range("H33").Value = 1085.4166
$ws.Range("I33").Value = 1167.2727
$ws.Range("J33").Value = 185
$ws.Range("K33").Value = 1167.2727
$ws.Range("L33").Value = 185
$ws.Range("M33").Value = -938.2727
$ws.Range("N33").Value = -643
$ws.Range("H40").Value = 11896.909
$ws.Range("I40").Value = 12586.7
$ws.Range("J40").Value = 4999
$ws.Range("K40").Value = 12586.7
$ws.Range("L40").Value = 4999
$ws.Range("M40").Value = -12411.7
$ws.Range("N40").Value = -5349
$ws.Range("H112").Value = 5284.3955
$ws.Range("J112").Value = 5697.615
$ws.Range("L112").Value = 17092.845
$ws.Range("N112").Value = -19308.845
$ws.Range("H129").Value = 1553.8
$ws.Range("I129").Value = 1553.8
$ws.Range("K129").Value = 4661.4
$ws.Range("M129").Value = 338.6000000000004
$ws.Range("H132").Value = 1822.56
$ws.Range("I132").Value = 1808.7347
$ws.Range("K132").Value = 5426.2041
$ws.Range("M132").Value = -2896.2041
$ws.Range("H137").Value = 3101.15
$ws.Range("I137").Value = 3565.625
$ws.Range("K137").Value = 10696.875
$ws.Range("M137").Value = -8146.875
$ws.Range("H138").Value = 5598.1055
$ws.Range("I138").Value = 2332.3333
$ws.Range("K138").Value = 6996.999899999999
$ws.Range("M138").Value = -1856.999899999999

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 10096.333
$ws.Range("I45").Value = 2256.5
$ws.Range("J45").Value = 12336.286
$ws.Range("K45").Value = 2256.5
$ws.Range("L45").Value = 12336.286
$ws.Range("M45").Value = -1879.5
$ws.Range("N45").Value = -13090.286
$ws.Range("H61").Value = 9682.210999999999
$ws.Range("I61").Value = 3042.4546
$ws.Range("K61").Value = 3042.4546
$ws.Range("M61").Value = -2830.4546
$ws.Range("H110").Value = 19610076
$ws.Range("I110").Value = 2276.5833
$ws.Range("K110").Value = 2276.5833
$ws.Range("M110").Value = -231.5832999999998
$ws.Range("H132").Value = 4911.943
$ws.Range("I132").Value = 2140.1304
$ws.Range("K132").Value = 6420.3912
$ws.Range("M132").Value = -3890.3912
$ws.Range("H136").Value = 9682.210999999999
$ws.Range("I136").Value = 3042.4546
$ws.Range("K136").Value = 9127.363799999999
$ws.Range("M136").Value = -6577.363799999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 345.13794
$ws.Range("I94").Value = 229.22223
$ws.Range("J94").Value = 534.8182
$ws.Range("K94").Value = 229.22223
$ws.Range("L94").Value = 534.8182
$ws.Range("M94").Value = 221.77777
$ws.Range("N94").Value = -1436.8182
$ws.Range("H134").Value = 8616.799999999999
$ws.Range("I134").Value = 4203
$ws.Range("K134").Value = 12609
$ws.Range("M134").Value = -10074

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 47619250
$ws.Range("I7").Value = 53.625
$ws.Range("J7").Value = 76923370
$ws.Range("K7").Value = 53.625
$ws.Range("L7").Value = 76923370
$ws.Range("M7").Value = 59.375
$ws.Range("N7").Value = -76923596
$ws.Range("H31").Value = 10827.087
$ws.Range("I31").Value = 4774.909
$ws.Range("K31").Value = 4774.909
$ws.Range("M31").Value = -4479.909
$ws.Range("H34").Value = 10827.087
$ws.Range("I34").Value = 4774.909
$ws.Range("K34").Value = 4774.909
$ws.Range("M34").Value = -4572.909
$ws.Range("H58").Value = 8599.259
$ws.Range("I58").Value = 3842.6667
$ws.Range("J58").Value = 9193.833000000001
$ws.Range("K58").Value = 3842.6667
$ws.Range("L58").Value = 9193.833000000001
$ws.Range("M58").Value = -3639.6667
$ws.Range("N58").Value = -9599.833000000001
$ws.Range("H99").Value = 4496.4443
$ws.Range("I99").Value = 3261.875
$ws.Range("J99").Value = 5484.1
$ws.Range("K99").Value = 3261.875
$ws.Range("L99").Value = 5484.1
$ws.Range("M99").Value = -1763.875
$ws.Range("N99").Value = -8480.1
$ws.Range("H126").Value = 4496.4443
$ws.Range("I126").Value = 3261.875
$ws.Range("J126").Value = 5484.1
$ws.Range("K126").Value = 9785.625
$ws.Range("L126").Value = 16452.3
$ws.Range("M126").Value = -7315.625
$ws.Range("N126").Value = -21392.3
$ws.Range("H132").Value = 7740.5454
$ws.Range("I132").Value = 4540.7
$ws.Range("K132").Value = 13622.1
$ws.Range("M132").Value = -11092.1
$ws.Range("H134").Value = 10767.258
$ws.Range("I134").Value = 10823.733
$ws.Range("J134").Value = 10714.3125
$ws.Range("K134").Value = 32471.199
$ws.Range("L134").Value = 32142.9375
$ws.Range("M134").Value = -29936.199
$ws.Range("N134").Value = -37212.9375
$ws.Range("H136").Value = 8599.259
$ws.Range("I136").Value = 3842.6667
$ws.Range("J136").Value = 9193.833000000001
$ws.Range("K136").Value = 11528.0001
$ws.Range("L136").Value = 27581.499
$ws.Range("M136").Value = -8978.000100000001
$ws.Range("N136").Value = -32681.499

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2877.1538
$ws.Range("I5").Value = 1334.3334
$ws.Range("J5").Value = 4199.5713
$ws.Range("K5").Value = 4003.0002
$ws.Range("L5").Value = 12598.7139
$ws.Range("M5").Value = -3891.0002
$ws.Range("N5").Value = -12822.7139
$ws.Range("H132").Value = 15558.909
$ws.Range("I132").Value = 9524.833000000001
$ws.Range("J132").Value = 22799.8
$ws.Range("K132").Value = 85723.497
$ws.Range("L132").Value = 205198.2
$ws.Range("M132").Value = -83193.497
$ws.Range("N132").Value = -210258.2
$ws.Range("H135").Value = 2877.1538
$ws.Range("I135").Value = 1334.3334
$ws.Range("J135").Value = 4199.5713
$ws.Range("K135").Value = 12009.0006
$ws.Range("L135").Value = 37796.14169999999
$ws.Range("M135").Value = -9474.000599999999
$ws.Range("N135").Value = -42866.14169999999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("H80").Value = 3780.0908
$ws.Range("I80").Value = 3434.3333
$ws.Range("J80").Value = 5336
$ws.Range("K80").Value = 3434.3333
$ws.Range("L80").Value = 5336
$ws.Range("M80").Value = -2436.3333
$ws.Range("N80").Value = -7332
$ws.Range("H83").Value = 3780.0908
$ws.Range("I83").Value = 3434.3333
$ws.Range("J83").Value = 5336
$ws.Range("K83").Value = 17171.6665
$ws.Range("L83").Value = 26680
$ws.Range("M83").Value = -12179.6665
$ws.Range("N83").Value = -36664
$ws.Range("H97").Value = 1167.5238
$ws.Range("I97").Value = 1001.4667
$ws.Range("K97").Value = 1001.4667
$ws.Range("M97").Value = -505.4666999999999
$ws.Range("H132").Value = 5705.0625
$ws.Range("I132").Value = 1364.6666
$ws.Range("K132").Value = 4093.9998
$ws.Range("M132").Value = -1563.9998
$ws.Range("M5").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4099.788
$ws.Range("I22").Value = 1846.4762
$ws.Range("J22").Value = 8043.0835
$ws.Range("K22").Value = 1846.4762
$ws.Range("L22").Value = 8043.0835
$ws.Range("M22").Value = -1551.4762
$ws.Range("N22").Value = -8633.083500000001
$ws.Range("H27").Value = 4099.788
$ws.Range("I27").Value = 1846.4762
$ws.Range("J27").Value = 8043.0835
$ws.Range("K27").Value = 1846.4762
$ws.Range("L27").Value = 8043.0835
$ws.Range("M27").Value = -1739.4762
$ws.Range("N27").Value = -8257.083500000001
$ws.Range("H46").Value = 1518.7307
$ws.Range("I46").Value = 826.8570999999999
$ws.Range("K46").Value = 826.8570999999999
$ws.Range("M46").Value = -638.8570999999999
$ws.Range("H68").Value = 142860590
$ws.Range("J68").Value = 4599.5
$ws.Range("L68").Value = 4599.5
$ws.Range("N68").Value = -6097.5
$ws.Range("H71").Value = 142860590
$ws.Range("J71").Value = 4599.5
$ws.Range("L71").Value = 22997.5
$ws.Range("N71").Value = -30485.5
$ws.Range("H122").Value = 3558.7222
$ws.Range("I122").Value = 2805.8262
$ws.Range("K122").Value = 8417.4786
$ws.Range("M122").Value = -5967.4786
$ws.Range("H132").Value = 6790.282
$ws.Range("I132").Value = 3583
$ws.Range("K132").Value = 10749
$ws.Range("M132").Value = -8219
$ws.Range("H136").Value = 10308.6455
$ws.Range("I136").Value = 5069.7144
$ws.Range("K136").Value = 15209.1432
$ws.Range("M136").Value = -12659.1432

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 14402381
$ws.Range("I122").Value = 20161566
$ws.Range("K122").Value = 60484698
$ws.Range("M122").Value = -60482248
$ws.Range("H132").Value = 44443
$ws.Range("I132").Value = 3331.6
$ws.Range("K132").Value = 9994.799999999999
$ws.Range("M132").Value = -7464.799999999999
$ws.Range("H136").Value = 34891.902
$ws.Range("I136").Value = 1888.2
$ws.Range("K136").Value = 5664.6
$ws.Range("M136").Value = -3114.6
